$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "only users who want to take on role of a " + "client" + " will
#    have an entry into this table" -> single merged run (same text).
#    Find across the run boundaries and replace with identical text;
#    the engine coalesces same-formatted runs into one on replace.
# ------------------------------------------------------------------
$rng = $d.Content
$target = "take on role of a client will have an entry into this table"
$found = $rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, $target, 1)
if (-not $found) {
    throw "Could not find paragraph 1 target text"
}

# ------------------------------------------------------------------
# 2) "Table 6: " (one bold run) -> "Table " / "5" / ": " (three bold
#    runs with identical rPr). Bookmarking the single character first
#    forces the run split to persist even though formatting matches
#    its neighbours; deleting the scratch bookmark afterwards leaves
#    the three runs intact without residual marks.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Table 6: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Table 6: ' heading"
}
$sixStart = $rng.Start + ("Table 6: ".IndexOf("6"))
$sixRange = $d.Range($sixStart, $sixStart + 1)
if ($sixRange.Text -ne "6") {
    throw "Offset arithmetic for 'Table 6' digit is wrong: got [$($sixRange.Text)]"
}
$d.Bookmarks.Add("zzzScratchSplit", $sixRange)
$digitRange = $d.Range($sixStart, $sixStart + 1)
$digitRange.Text = "5"
$d.Bookmarks("zzzScratchSplit").Delete()

# ------------------------------------------------------------------
# 3) "Table 7: messages" (one bold run) -> "Table " / "7" / ": messages"
#    (three bold runs) with the "_GoBack" bookmark collapsed between
#    the "7" run and the ": messages" run. Adding a bookmark named
#    "_GoBack" automatically relocates the existing one (Word keeps a
#    single "_GoBack" bookmark), which also removes it from its old
#    location at the end of the document.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Table 7: messages", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Table 7: messages' heading"
}
$sevenStart = $rng.Start + ("Table 7: messages".IndexOf("7"))
$sevenRange = $d.Range($sevenStart, $sevenStart + 1)
if ($sevenRange.Text -ne "7") {
    throw "Offset arithmetic for 'Table 7' digit is wrong: got [$($sevenRange.Text)]"
}
$d.Bookmarks.Add("zzzScratchSplit2", $sevenRange)
$d.Bookmarks("zzzScratchSplit2").Delete()
$goBackRange = $d.Range($sevenStart + 1, $sevenStart + 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# ------------------------------------------------------------------
# 4) "recipient" + "_id" -> single merged run "recipient_id" (same text).
# ------------------------------------------------------------------
$rng = $d.Content
$target = "recipient_id"
$found = $rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, $target, 1)
if (-not $found) {
    throw "Could not find 'recipient_id' text"
}

Write-Output "done"
